$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "53.800.77"
$ws.Range("E2").Value = "  -4.78%  "

Set-TextValue "D3" "2.237.68"
$ws.Range("E3").Value = "  -6.02%  "

$ws.Range("E4").Value = "  +0.03%  "

Set-TextValue "D5" "487.45"
$ws.Range("E5").Value = "  -3.00%  "

Set-TextValue "D6" "126.38"
$ws.Range("E6").Value = "  -3.06%  "

$ws.Range("E7").Value = "  +0.22%  "

Set-TextValue "D8" "0.520"
$ws.Range("E8").Value = "  -4.69%  "

Set-TextValue "D9" "2.250.01"
$ws.Range("E9").Value = "  -5.84%  "

$ws.Range("E10").Value = "  -7.42%  "

$ws.Range("E11").Value = "  -1.51%  "

Set-TextValue "D12" "4.73"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("E13").Value = "  -4.18%  "

Set-TextValue "D14" "2.636.87"
$ws.Range("E14").Value = "  -5.95%  "

Set-TextValue "D15" "21.09"
$ws.Range("E15").Value = "  -2.85%  "

Set-TextValue "D16" "53.713.20"
$ws.Range("E16").Value = "  -4.86%  "

$ws.Range("E17").Value = "  -3.64%  "

Set-TextValue "D18" "2.238.79"
$ws.Range("E18").Value = "  -6.46%  "

$ws.Range("E19").Value = "  -1.55%  "

Set-TextValue "D20" "9.58"
$ws.Range("E20").Value = "  -4.96%  "

Set-TextValue "D21" "299.08"
$ws.Range("E21").Value = "  -2.88%  "

Set-TextValue "D22" "6.12"
$ws.Range("E22").Value = "  -2.43%  "

Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  +0.12%  "

Set-TextValue "D24" "63.64"
$ws.Range("E24").Value = "  -2.47%  "

Set-TextValue "D25" "0.999"
$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("E26").Value = "  -1.27%  "

$ws.Range("E27").Value = "  -3.55%  "

Set-TextValue "D28" "7.03"
$ws.Range("E28").Value = "  -4.04%  "

Set-TextValue "D29" "169.10"
$ws.Range("E29").Value = "  -2.10%  "

Set-TextValue "D30" "0.0₃0688"
$ws.Range("E30").Value = "  -3.95%  "

$ws.Range("E31").Value = "  -2.90%  "

Set-TextValue "D32" "0.998"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("E33").Value = "  -0.51%  "

Set-TextValue "D34" "0.997"
$ws.Range("E34").Value = "  +0.09%  "

Set-TextValue "D35" "1.05"
$ws.Range("E35").Value = "  -3.08%  "

Set-TextValue "D36" "17.45"
$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("E37").Value = "  -1.46%  "

Set-TextValue "D38" "0.841"
$ws.Range("E38").Value = "  +5.76%  "

Set-TextValue "D39" "3.57"
$ws.Range("E39").Value = "  -5.64%  "

Set-TextValue "D40" "35.65"
$ws.Range("E40").Value = "  -1.19%  "

Set-TextValue "D41" "0.366"
$ws.Range("E41").Value = "  -0.82%  "

Set-TextValue "D42" "1.36"
$ws.Range("E42").Value = "  -2.72%  "

$ws.Range("E43").Value = "  -2.02%  "

Set-TextValue "D44" "122.69"
$ws.Range("E44").Value = "  -6.39%  "

Set-TextValue "D45" "4.66"
$ws.Range("E45").Value = "  -2.34%  "

Set-TextValue "D46" "0.0879"
$ws.Range("E46").Value = "  -3.28%  "

$ws.Range("E47").Value = "  -5.54%  "

Set-TextValue "D48" "235.68"
$ws.Range("E48").Value = "  -2.61%  "

Set-TextValue "D49" "0.0470"
$ws.Range("E49").Value = "  -2.98%  "

Set-TextValue "D50" "0.0202"
$ws.Range("E50").Value = "  -3.64%  "

$ws.Range("E51").Value = "  -4.65%  "
